# "Template fix and svg file"
#
# Changes applied to the workbook:
#  1. Rename the "pdx_models" sheet to "pdx_model" (template fix).
#  2. Make "pdx_model" the active/selected sheet (instead of "cell_model"),
#     which flips the `tabSelected` flag between the two sheets and moves
#     the workbook's `activeTab` pointer.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("pdx_models")
$ws.Name = "pdx_model"

# Activating the renamed sheet makes it the workbook's active tab (tabSelected="1"
# on its sheetView) and clears the flag from whichever sheet previously held it
# ("cell_model").
$ws.Activate()
